$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its literal text formatting (values like "1.000",
# "26.675.28", "332.49" must not be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.675.28"
$ws.Range("E2").Value = "  +7.32%  "

$ws.Range("D3").Value = "1.741.00"
$ws.Range("E3").Value = "  +5.04%  "

$ws.Range("D4").Value = "1.003"

$ws.Range("D5").Value = "332.49"
$ws.Range("E5").Value = "  +6.85%  "

$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").Value = "0.3745"
$ws.Range("E7").Value = "  +3.56%  "

$ws.Range("D8").Value = "49.10"
$ws.Range("E8").Value = "  +4.36%  "

$ws.Range("D9").Value = "0.3418"
$ws.Range("E9").Value = "  +5.32%  "

$ws.Range("D10").Value = "1.199"
$ws.Range("E10").Value = "  +6.22%  "

$ws.Range("D11").Value = "0.07496"
$ws.Range("E11").Value = "  +6.71%  "

$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").Value = "6.453"
$ws.Range("E13").Value = "  +7.62%  "

$ws.Range("D14").Value = "20.39"
$ws.Range("E14").Value = "  +5.02%  "

$ws.Range("D15").Value = "7.142"
$ws.Range("E15").Value = "  +8.67%  "

$ws.Range("D16").Value = "1.739.22"
$ws.Range("E16").Value = "  +4.44%  "

$ws.Range("D17").Value = "0.00001090"
$ws.Range("E17").Value = "  +4.87%  "

$ws.Range("D18").Value = "0.06690"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").Value = "83.68"
$ws.Range("E19").Value = "  +6.25%  "

$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").Value = "16.74"
$ws.Range("E21").Value = "  +7.29%  "

$ws.Range("D22").Value = "6.199"
$ws.Range("E22").Value = "  +5.46%  "

$ws.Range("D23").Value = "13.11"
$ws.Range("E23").Value = "  +4.56%  "

$ws.Range("D24").Value = "26.666.63"
$ws.Range("E24").Value = "  +7.40%  "

$ws.Range("D25").Value = "2.449"
$ws.Range("E25").Value = "  +1.13%  "

$ws.Range("D26").Value = "2.487"
$ws.Range("E26").Value = "  +3.72%  "

$ws.Range("D27").Value = "1.413"
$ws.Range("E27").Value = "  +18.93%  "

$ws.Range("D28").Value = "153.35"
$ws.Range("E28").Value = "  +4.41%  "

$ws.Range("D29").Value = "19.67"
$ws.Range("E29").Value = "  +6.42%  "

$ws.Range("D30").Value = "1.936.06"
$ws.Range("E30").Value = "  +4.61%  "

$ws.Range("D31").Value = "131.84"
$ws.Range("E31").Value = "  +5.85%  "

$ws.Range("D32").Value = "4.139"
$ws.Range("E32").Value = "  +1.04%  "

$ws.Range("D33").Value = "6.105"
$ws.Range("E33").Value = "  +7.35%  "

$ws.Range("D34").Value = "0.08604"
$ws.Range("E34").Value = "  +2.11%  "

$ws.Range("D35").Value = "1.711"
$ws.Range("E35").Value = "  +2.89%  "

$ws.Range("D36").Value = "13.02"
$ws.Range("E36").Value = "  +6.86%  "

$ws.Range("D37").Value = "5.453"
$ws.Range("E37").Value = "  +6.69%  "

$ws.Range("D38").Value = "0.02370"
$ws.Range("E38").Value = "  +5.97%  "

$ws.Range("D39").Value = "0.06334"
$ws.Range("E39").Value = "  +5.54%  "

$ws.Range("D40").Value = "0.2186"
$ws.Range("E40").Value = "  +6.77%  "

$ws.Range("D41").Value = "8.653"
$ws.Range("E41").Value = "  +4.86%  "

$ws.Range("D42").Value = "1.236"
$ws.Range("E42").Value = "  -3.63%  "

$ws.Range("E43").Value = "  +6.63%  "

$ws.Range("D44").Value = "14.31"
$ws.Range("E44").Value = "  +13.08%  "

$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "3.911"
$ws.Range("E46").Value = "  +3.82%  "

$ws.Range("D47").Value = "0.6081"
$ws.Range("E47").Value = "  +9.37%  "

$ws.Range("D48").Value = "129.57"
$ws.Range("E48").Value = "  +4.22%  "

$ws.Range("D49").Value = "2.074"
$ws.Range("E49").Value = "  +7.46%  "

$ws.Range("D50").Value = "0.07299"
$ws.Range("E50").Value = "  +4.65%  "

$ws.Range("D51").Value = "77.82"
$ws.Range("E51").Value = "  +4.81%  "
